$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) IC2 (row 4): footprint fix from TXS0108EPWR / TSSOP-20 to the
#    correct quad level-shifter part + its TSSOP-24 footprint.
# -----------------------------------------------------------------
$ws.Range("A4").Value = "74LVC4245APWR"
$ws.Range("C4").Value = "TSSOP-24_4.4x7.8x0.65P"
$ws.Range("D4").Value = "C7859"

# -----------------------------------------------------------------
# 2) C2 (row 6): simplify footprint name.
# -----------------------------------------------------------------
$ws.Range("C6").Value = "3528"

# -----------------------------------------------------------------
# 3) Insert a new BOM row (new row 8) for the single-gate transceiver
#    U2, pushing everything below it down by one row.
# -----------------------------------------------------------------
$ws.Range("A8:D8").Insert(-4121)  # xlShiftDown

# Copy the formatting (fonts/fills/borders/number format) from the row
# that used to be row 8 (now shifted down to row 9) so the new row
# matches the rest of the table exactly.
$ws.Range("A9:D9").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows("8:8").RowHeight = $ws.Rows("9:9").RowHeight

$ws.Range("A8").Value = "74LVC2T45DCURG4"
$ws.Range("B8").Value = "U2"
$ws.Range("C8").Value = "VSSOP-8_2.0x2.3x0.5P"
$ws.Range("D8").Value = "C139382"

# -----------------------------------------------------------------
# 4) R1 (now row 13): reduce LED brightness by raising resistor value.
# -----------------------------------------------------------------
$ws.Range("A13").Value = "82Ω ±1% 1/8W"
$ws.Range("D13").Value = "C17841"

# -----------------------------------------------------------------
# 5) R2 R3 R4 (now row 14): reduce LED brightness by raising resistor value.
# -----------------------------------------------------------------
$ws.Range("A14").Value = "330Ω ±1% 1/8W"
$ws.Range("D14").Value = "C17630"

# -----------------------------------------------------------------
# 6) 100NF 50V row (now row 15): more capacitors share this BOM line.
# -----------------------------------------------------------------
$ws.Range("B15").Value = "C4 C5 C6 C7 C8 C9"
